# "respuestas salida.xlsx" update
# The author expanded several short chatbot "RESPUESTA" answers into longer,
# more polite/complete phrasings while keeping the same "CLASE" category for
# each row. Six of the fifteen data rows (A2:A15) get their answer text
# replaced; everything else (headers, other rows, categories) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the six changed "RESPUESTA" cells, in the same order the strings
# were (re)written to the shared-strings table.
$ws.Range("A4").Value  = "Hola cómo está, en qué puedo ayudarlo"
$ws.Range("A5").Value  = "Buenos días, en qué puedo ayudarlo"
$ws.Range("A13").Value = "Me llamo Robotin, estoy a sus ordenes"
$ws.Range("A14").Value = "Mi nombre es Robotin y hoy les estaré ayudando"
$ws.Range("A15").Value = "Robotin es mi nombre, cómo puedo ayudarle"
$ws.Range("A12").Value = "Mi nombre Robotin, estoy acá para ayudarle"

# Column A ("RESPUESTA") holds the longer text now, so re-fit its width to
# the new content, same as double-clicking the column border in Excel.
$ws.Columns("A:A").ColumnWidth = 39.6

# The author's cursor ended up on D14 when the file was last saved.
$ws.Range("D14").Select() | Out-Null
